$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.255537033081055
$ws.Range("B1").Value = 2.382030963897705
$ws.Range("D1").Value = 1.395228147506714
$ws.Range("E1").Value = 0.8727909326553345
